# Fatura satır tablosu güncellemesi:
#  - 2. satırdaki ürün adı, birim fiyat ve KDV tutarı değiştirildi
#  - 3. satıra yeni bir ürün satırı eklendi (önceden boştu)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update product name, unit price and KDV amount
$ws.Range("A2").Value = "a (b)"
$ws.Range("D2").Value = 150
$ws.Range("H2").Value = 30

# Row 3: fill in the previously-empty row with a new line item
$ws.Range("A3").Value = "d (e)"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "C62"
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 20
